$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting existing rows 3-18 down to 4-19.
# (Excel's default "shift cells down" insert, inheriting formatting from
# the row above - matches the style already applied to column D dates.)
$ws.Range("A3:R3").Insert()

# Populate the newly inserted row 3 with the new weekly price entry.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44831
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112012
$ws.Range("G3").Value = "Espinaca"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 7500
$ws.Range("M3").Value = 7250
$ws.Range("N3").Value = "$/cuna 10 kilos"
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 725
$ws.Range("Q3").Value = 10
$ws.Range("R3").Value = "Hortaliza"
